# Inserts a new weekly price record for "Betarraga" at Feria Lagunitas de
# Puerto Montt. The new observation is dated 44627 and slots in chronologically
# where row 197 used to be, pushing the existing rows 197-261 down to 198-262
# (dimension grows from A1:R261 to A1:R262).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 197; this shifts rows 197..261 down to 198..262
# and extends the sheet dimension automatically.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new price record.
$ws.Cells.Item(197, 1).Value = 4
$ws.Cells.Item(197, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(197, 3).Value = "Los Lagos"
$ws.Cells.Item(197, 4).Value = 44627
$ws.Cells.Item(197, 5).Value = 10
$ws.Cells.Item(197, 6).Value = 100114014
$ws.Cells.Item(197, 7).Value = "Betarraga"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 500
$ws.Cells.Item(197, 11).Value = 800
$ws.Cells.Item(197, 12).Value = 1000
$ws.Cells.Item(197, 13).Value = 900
$ws.Cells.Item(197, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(197, 15).Value = "Región del Maule"
$ws.Cells.Item(197, 16).Value = 180
$ws.Cells.Item(197, 17).Value = 5
$ws.Cells.Item(197, 18).Value = "Hortaliza"
